$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "863268050609451"
$ws.Range("A3").NumberFormat = "@"
$ws.Range("A3").Value = "863268050609683"

$ws.Range("A4").Select()
